$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.538.27"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").Value = "2.515.76"
$ws.Range("E3").Value = "  -4.41%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "583.65"
$ws.Range("E5").Value = "  -1.84%  "
$ws.Range("D6").Value = "173.90"
$ws.Range("E6").Value = "  +3.44%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "0.524"
$ws.Range("E8").Value = "  -1.85%  "
$ws.Range("D9").Value = "2.514.93"
$ws.Range("E9").Value = "  -4.40%  "
$ws.Range("E10").Value = "  -0.43%  "
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("D12").Value = "0.351"
$ws.Range("E12").Value = "  -3.55%  "
$ws.Range("E13").Value = "  -2.35%  "
$ws.Range("D14").Value = "26.65"
$ws.Range("E14").Value = "  -3.70%  "
$ws.Range("D15").Value = "2.962.05"
$ws.Range("E15").Value = "  -4.76%  "
$ws.Range("D16").Value = "0.0000177"
$ws.Range("E16").Value = "  -2.91%  "
$ws.Range("D17").Value = "66.235.46"
$ws.Range("E17").Value = "  -1.75%  "
$ws.Range("D18").Value = "2.510.15"
$ws.Range("E18").Value = "  -4.39%  "
$ws.Range("D19").Value = "7.74"
$ws.Range("E19").Value = "  -3.74%  "
$ws.Range("D20").Value = "11.28"
$ws.Range("E20").Value = "  -5.88%  "
$ws.Range("D21").Value = "348.62"
$ws.Range("E21").Value = "  -2.40%  "
$ws.Range("D22").Value = "4.20"
$ws.Range("E22").Value = "  -2.57%  "
$ws.Range("D23").Value = "4.61"
$ws.Range("E23").Value = "  -1.17%  "
$ws.Range("D24").Value = "1.98"
$ws.Range("E24").Value = "  +2.08%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").Value = "69.70"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "9.95"
$ws.Range("E27").Value = "  -3.64%  "
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").Value = "2.642.79"
$ws.Range("E29").Value = "  -4.52%  "
$ws.Range("D30").Value = "0.0₃0978"
$ws.Range("E30").Value = "  -3.08%  "
$ws.Range("D31").Value = "528.63"
$ws.Range("E31").Value = "  -3.51%  "
$ws.Range("D32").Value = "8.13"
$ws.Range("E32").Value = "  +2.44%  "
$ws.Range("E33").Value = "  -2.45%  "
$ws.Range("D34").Value = "1.84"
$ws.Range("E34").Value = "  -3.09%  "
$ws.Range("E35").Value = "  -3.76%  "
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("E37").Value = "  -2.75%  "
$ws.Range("D38").Value = "155.81"
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("D39").Value = "18.60"
$ws.Range("E39").Value = "  -2.25%  "
$ws.Range("D40").Value = "18.35"
$ws.Range("E40").Value = "  +0.31%  "
$ws.Range("D41").Value = "0.355"
$ws.Range("E41").Value = "  -2.71%  "
$ws.Range("E42").Value = "  -1.88%  "
$ws.Range("E43").Value = "  -2.25%  "
$ws.Range("D46").Value = "39.64"
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("D47").Value = "147.98"
$ws.Range("E47").Value = "  -3.33%  "
$ws.Range("D48").Value = "0.558"
$ws.Range("E48").Value = "  -3.69%  "
$ws.Range("D51").Value = "0.0₆0272"
$ws.Range("E51").Value = "  -8.96%  "

# Row swaps (coin rows exchanged with neighbor)
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "2.54"
$ws.Range("E44").Value = "  +4.47%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").Value = "3.67"
$ws.Range("E49").Value = "  -3.32%  "
$ws.Range("B50").Value = "Optimism"
$ws.Range("C50").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D50").Value = "1.73"
$ws.Range("E50").Value = "  +1.66%  "

Write-Output "Applied cryptos update"
